$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / info block text updates -------------------------------------
$ws.Range("A1").Value = "Evento Bienestar"
$ws.Range("A6").Value = "Instructor: enzy zulay angarita bermudez"
$ws.Range("A7").Value = "Clase Formacion: Algoritmia"
$ws.Range("A8").Value = "Fecha: 2024-10-10 07:46:15"

# --- Existing attendee row: hours-absent value + color update -------------
$ws.Range("H12").Value = 1
$ws.Range("H12").Interior.Color = 39423

# --- New attendee row ------------------------------------------------------
$ws.Range("A13").Value = 3

# B13 / D13 look like numbers ("0980987651", "3128765423") but must be kept
# as text (so the leading zero in the document number survives) - force the
# cell to Text format before assigning, then strip the format back off so
# the cells don't carry a stray number format.
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "0980987651"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3128765423"
$ws.Range("B13:D13").ClearFormats()

$ws.Range("C13").Value = "Miguel Alexander Toloza"
$ws.Range("E13").Value = "miguel@gmail.com"
$ws.Range("F13").Value = "Masculino"
$ws.Range("G13").Value = "Santander - Bucaramanga - La Esperanza"
$ws.Range("H13").Value = 1
$ws.Range("H13").Interior.Color = 39423

"done"
